$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new support value and itemset label
$ws.Range("A2").Value = 0.1687841891653432
$ws.Range("B2").Value = "frozenset({'Mpx 2 10W 30 Sl 800 mL'})"

# Remove rows 3 and 4 entirely, shrinking the used range to A1:B2
$ws.Range("A3:B4").Delete()
